$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.398.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "
# Row 3
$ws.Range("D3").Value = "'1.869.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "
# Row 5
$ws.Range("D5").Value = "'243.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
# Row 6
$ws.Range("D6").Value = "'0.7039"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.97%  "
# Row 7
$ws.Range("E7").Value = "  -0.12%  "
# Row 8
$ws.Range("D8").Value = "'0.07924"
$ws.Range("D8").Style = "Normal"
# Row 9
$ws.Range("D9").Value = "'0.3132"
$ws.Range("D9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = "'24.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.77%  "
# Row 11
$ws.Range("D11").Value = "'0.07838"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.52%  "
# Row 12
$ws.Range("D12").Value = "'1.906.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "
# Row 13
$ws.Range("D13").Value = "'93.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.87%  "
# Row 14
$ws.Range("D14").Value = "'5.171"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.99%  "
# Row 15
$ws.Range("D15").Value = "'0.7009"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.47%  "
# Row 16
$ws.Range("D16").Value = "'6.522"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.20%  "
# Row 17
$ws.Range("D17").Value = "'0.000008403"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.05%  "
# Row 18
$ws.Range("D18").Value = "'29.486.92"
$ws.Range("D18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = "'252.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.88%  "
# Row 20
$ws.Range("D20").Value = "'2.147.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "
# Row 21
$ws.Range("E21").Value = "  -1.02%  "
# Row 22
$ws.Range("E22").Value = "  -0.17%  "
# Row 23
$ws.Range("D23").Value = "'7.667"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.22%  "
# Row 24
$ws.Range("E24").Value = "  -0.21%  "
# Row 25
$ws.Range("D25").Value = "'0.1553"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.76%  "
# Row 26
$ws.Range("D26").Value = "'9.015"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
# Row 27
$ws.Range("E27").Value = "  -0.52%  "
# Row 28
$ws.Range("E28").Value = "  +1.74%  "
# Row 29
$ws.Range("D29").Value = "'1.507"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.35%  "
# Row 30
$ws.Range("E30").Value = "  -2.02%  "
# Row 31
$ws.Range("D31").Value = "'4.261"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.94%  "
# Row 32
$ws.Range("E32").Value = "  +2.26%  "
# Row 33
$ws.Range("D33").Value = "'0.05266"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.65%  "
# Row 34
$ws.Range("D34").Value = "'1.897"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.02%  "
# Row 35
$ws.Range("D35").Value = "'1.180"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.40%  "
# Row 36
$ws.Range("D36").Value = "'0.7506"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "
# Row 37
$ws.Range("D37").Value = "'2.711"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.35%  "
# Row 38
$ws.Range("D38").Value = "'0.01878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
# Row 39
$ws.Range("D39").Value = "'1.272.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "
# Row 40
$ws.Range("E40").Value = "  +0.56%  "
# Row 41
$ws.Range("D41").Value = "'0.8925"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.76%  "
# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'109.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.24%  "
# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.041"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.95%  "
# Row 44
$ws.Range("D44").Value = "'70.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.29%  "
# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "'2.041.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "
# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000126"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.04%  "
# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.805"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.67%  "
# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.631"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "
# Row 50
$ws.Range("D50").Value = "'0.5184"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.84%  "
# Row 51
$ws.Range("E51").Value = "  -0.82%  "
